# Removed Test Case Inter-Dependency
#
# The "productname" value (row 1, column B) is shared between the
# ProductLoanInput and ProductLoanOutput sheets via the same shared
# string; update it on both sheets so the shared string itself is
# edited (rather than orphaned) and both sheets stay in sync.
#
# The "shortname" value (row 2, column B) on ProductLoanInput switches
# from the numeric literal 4301 to a new independent text value "430w"
# so the short name test data is no longer derived from/tied to the
# product name/id used elsewhere.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4301-MS-EI-DB-SAR-REC-RNI-FEE+INT-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DL-FIFR-1-MD-TR-1st"

$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

$wsInput.Range("B2").Value = "430w"

# Make ProductLoanInput the active/selected tab (it was ProductLoanOutput
# before), with B3 as the selected cell. ProductLoanOutput's own
# selection (B1) is already at its original/default value so it does
# not need to be touched.
$wsInput.Activate()
$wsInput.Range("B3").Select()
